$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-15 (columns A-F) reflecting the rolling-window update
$data = @(
    @(45944, 5564.50698273871, 5684.6799371782, 3620, 8104.078361, 191.843804809979),
    @(45945, 5816.11417805593, 5915.65361082697, 3620, 8539.47947, 209.12578761546),
    @(45946, 6218.61815962666, 6241.36067077709, 3620, 9246.001595, 235.364337756268),
    @(45947, 6118.17196890358, 5630.11098002927, 3620, 9082.815358, 207.281432046904),
    @(45948, 1742.27770790123, 4004.92596044467, 3620, 4473.853115, 129.854223647643),
    @(45949, 1636.94065696827, 4042.60018843985, 3620, 4358.692014, 131.014647727982),
    @(45950, 6392.95297294923, 6558.77377325878, 3620, 9928.522437, 269.764301554564),
    @(45951, 6392.95297294923, 6679.26230254023, 3620, 9928.522437, 274.784656941291),
    @(45952, 6392.95297294923, 6570.9056564803, 3620, 9928.522437, 270.269796688794),
    @(45953, 6392.95297294923, 6287.82268472799, 3620, 9928.522437, 258.474672865781),
    @(45954, 6392.95297294923, 5639.10881235012, 3620, 9928.522437, 231.44492818337),
    @(45955, 2044.14689306201, 4024.22706818106, 3620, 5157.782237, 146.57760050496),
    @(45956, 1933.9611505005, 3930.73269626633, 3620, 5037.335348, 142.25445390691),
    @(45957, 6757.64823424925, 6304.09430704083, 3620, 10717.728041, 276.84058807465)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
}
